$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header cell F1, using the same format as the other header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "19-jun"

# Fill column F (rows 2-25) with placeholder "-" values for the new day
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = "-"
}
